$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1717.0416
$ws.Range("J70").Value = 1743.0952
$ws.Range("L70").Value = 5229.2856
$ws.Range("N70").Value = -5769.2856
$ws.Range("H73").Value = 1717.0416
$ws.Range("J73").Value = 1743.0952
$ws.Range("L73").Value = 5229.2856
$ws.Range("N73").Value = -7101.2856
$ws.Range("H116").Value = 5769.067
$ws.Range("I116").Value = 8564.823
$ws.Range("K116").Value = 8564.823
$ws.Range("M116").Value = -5122.823
$ws.Range("H132").Value = 184428.44
$ws.Range("I132").Value = 3029.6829
$ws.Range("J132").Value = 715667.6
$ws.Range("K132").Value = 9089.048699999999
$ws.Range("L132").Value = 2147002.8
$ws.Range("M132").Value = -6559.048699999999
$ws.Range("N132").Value = -2152062.8
$ws.Range("H135").Value = 10205435
$ws.Range("I135").Value = 359.51852
$ws.Range("J135").Value = 22729846
$ws.Range("K135").Value = 3235.66668
$ws.Range("L135").Value = 204568614
$ws.Range("M135").Value = -700.6666800000003
$ws.Range("N135").Value = -204573684
$ws.Range("H137").Value = 26746.41
$ws.Range("I137").Value = 50879.35
$ws.Range("J137").Value = 1343.3158
$ws.Range("K137").Value = 152638.05
$ws.Range("L137").Value = 4029.9474
$ws.Range("M137").Value = -150088.05
$ws.Range("N137").Value = -9129.947400000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 60014
$ws.Range("I34").Value = 45000
$ws.Range("J34").Value = 75028
$ws.Range("K34").Value = 45000
$ws.Range("L34").Value = 75028
$ws.Range("M34").Value = -44729
$ws.Range("N34").Value = -75570
$ws.Range("H62").Value = 87686.75
$ws.Range("J62").Value = 87686.75
$ws.Range("L62").Value = 87686.75
$ws.Range("N62").Value = -88934.75
$ws.Range("H65").Value = 87686.75
$ws.Range("J65").Value = 87686.75
$ws.Range("L65").Value = 263060.25
$ws.Range("N65").Value = -269300.25
$ws.Range("H97").Value = 782.7368
$ws.Range("I97").Value = 417.65
$ws.Range("J97").Value = 1188.3889
$ws.Range("K97").Value = 417.65
$ws.Range("L97").Value = 1188.3889
$ws.Range("M97").Value = 78.35000000000002
$ws.Range("N97").Value = -2180.3889

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 293164.12
$ws.Range("I86").Value = 1496.6
$ws.Range("J86").Value = 1751501.8
$ws.Range("K86").Value = 1496.6
$ws.Range("L86").Value = 1751501.8
$ws.Range("M86").Value = -373.5999999999999
$ws.Range("N86").Value = -1753747.8
$ws.Range("H89").Value = 293164.12
$ws.Range("I89").Value = 1496.6
$ws.Range("J89").Value = 1751501.8
$ws.Range("K89").Value = 7483
$ws.Range("L89").Value = 8757509
$ws.Range("M89").Value = -1867
$ws.Range("N89").Value = -8768741
$ws.Range("H94").Value = 1453.7097
$ws.Range("I94").Value = 364.47058
$ws.Range("J94").Value = 2776.3572
$ws.Range("K94").Value = 364.47058
$ws.Range("L94").Value = 2776.3572
$ws.Range("M94").Value = 86.52942000000002
$ws.Range("N94").Value = -3678.3572
$ws.Range("H99").Value = 1051.1428
$ws.Range("I99").Value = 1093
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 1093
$ws.Range("L99").Value = 800
$ws.Range("M99").Value = 405
$ws.Range("N99").Value = -3796

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12337.857
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 12337.857
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 12337.857
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -12927.857
$ws.Range("H34").Value = 12337.857
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 12337.857
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12337.857
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -12741.857
$ws.Range("H63").Value = 52708.4
$ws.Range("J63").Value = 52708.4
$ws.Range("L63").Value = 52708.4
$ws.Range("N63").Value = -54080.4
$ws.Range("H66").Value = 52708.4
$ws.Range("J66").Value = 52708.4
$ws.Range("L66").Value = 158125.2
$ws.Range("N66").Value = -164989.2

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 803.5454999999999
$ws.Range("I5").Value = 654.5714
$ws.Range("J5").Value = 913.3158
$ws.Range("K5").Value = 1963.7142
$ws.Range("L5").Value = 2739.9474
$ws.Range("M5").Value = -1851.7142
$ws.Range("N5").Value = -2963.9474
$ws.Range("H60").Value = 212.8
$ws.Range("I60").Value = 212.8
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 638.4000000000001
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -387.4000000000001
$ws.Range("N60").Value = $null
$ws.Range("H122").Value = 527.8333
$ws.Range("I122").Value = 305.75
$ws.Range("J122").Value = 749.9167
$ws.Range("K122").Value = 2751.75
$ws.Range("L122").Value = 6749.2503
$ws.Range("M122").Value = -301.75
$ws.Range("N122").Value = -11649.2503
$ws.Range("H135").Value = 803.5454999999999
$ws.Range("I135").Value = 654.5714
$ws.Range("J135").Value = 913.3158
$ws.Range("K135").Value = 5891.1426
$ws.Range("L135").Value = 8219.842199999999
$ws.Range("M135").Value = -3356.1426
$ws.Range("N135").Value = -13289.8422

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5782.143
$ws.Range("I80").Value = 2454.3635
$ws.Range("J80").Value = 7935.4116
$ws.Range("K80").Value = 2454.3635
$ws.Range("L80").Value = 7935.4116
$ws.Range("M80").Value = -1456.3635
$ws.Range("N80").Value = -9931.411599999999
$ws.Range("H83").Value = 5782.143
$ws.Range("I83").Value = 2454.3635
$ws.Range("J83").Value = 7935.4116
$ws.Range("K83").Value = 12271.8175
$ws.Range("L83").Value = 39677.058
$ws.Range("M83").Value = -7279.817499999999
$ws.Range("N83").Value = -49661.058

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3960.8
$ws.Range("I61").Value = 4601.3335
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 4601.3335
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -4399.3335
$ws.Range("N61").Value = -3404
$ws.Range("H100").Value = 23803.936
$ws.Range("I100").Value = 49256.24
$ws.Range("J100").Value = 2424
$ws.Range("K100").Value = 49256.24
$ws.Range("L100").Value = 2424
$ws.Range("M100").Value = -48715.24
$ws.Range("N100").Value = -3506
$ws.Range("H113").Value = 3960.8
$ws.Range("I113").Value = 4601.3335
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 4601.3335
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2431.3335
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 50469.5
$ws.Range("I132").Value = 67007.13
$ws.Range("J132").Value = 3863.4546
$ws.Range("K132").Value = 201021.39
$ws.Range("L132").Value = 11590.3638
$ws.Range("M132").Value = -198491.39
$ws.Range("N132").Value = -16650.3638

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 250.35
$ws.Range("I113").Value = 250.44444
$ws.Range("K113").Value = 751.33332
$ws.Range("M113").Value = 1418.66668
